# Applies crypto price/volume updates per commit "Updated cryptos list on Wed Sep 13 15:20:56 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.292.53"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").Value = "1.608.99"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.02"
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("E8").Value = "  -0.23%  "
$ws.Range("E9").Value = "  -0.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.40"
$ws.Range("E10").Value = "  +1.80%  "
$ws.Range("E11").Value = "  -0.07%  "
$ws.Range("D12").Value = "1.832.31"
$ws.Range("E12").Value = "  +0.18%  "
$ws.Range("D13").Value = "1.611.87"
$ws.Range("E13").Value = "  +0.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.05"
$ws.Range("E14").Value = "  +0.75%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.515"
$ws.Range("E15").Value = "  +0.87%  "
$ws.Range("D16").Value = "26.302.14"
$ws.Range("E16").Value = "  +0.40%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.78"
$ws.Range("E17").Value = "  +1.90%  "
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("E19").Value = "  -0.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "203.46"
$ws.Range("E20").Value = "  +2.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.30"
$ws.Range("E21").Value = "  +1.36%  "
$ws.Range("E22").Value = "  -0.90%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("E24").Value = "  +8.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.33"
$ws.Range("E25").Value = "  +1.08%  "
$ws.Range("E26").Value = "  -0.16%  "
$ws.Range("E27").Value = "  -4.70%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.24"
$ws.Range("E28").Value = "  +0.49%  "
$ws.Range("E29").Value = "  +1.43%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0489"
$ws.Range("E30").Value = "  +3.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.17"
$ws.Range("E31").Value = "  -0.19%  "
$ws.Range("E32").Value = "  +1.94%  "
$ws.Range("E33").Value = "  -2.16%  "
$ws.Range("E34").Value = "  +3.31%  "
$ws.Range("E35").Value = "  +0.29%  "
$ws.Range("D36").Value = "1.159.62"
$ws.Range("E36").Value = "  +4.69%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0167"
$ws.Range("E37").Value = "  +9.93%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.798"
$ws.Range("E38").Value = "  +1.60%  "
$ws.Range("B39").Value = "PaxDollar"
$ws.Range("C39").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  -0.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.32"
$ws.Range("E40").Value = "  -0.33%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.502"
$ws.Range("E41").Value = "  +0.70%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.786"
$ws.Range("E42").Value = "  +1.42%  "
$ws.Range("E43").Value = "  +2.88%  "
$ws.Range("D44").Value = "1.746.02"
$ws.Range("E44").Value = "  +0.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.71"
$ws.Range("E45").Value = "  -0.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.54"
$ws.Range("E46").Value = "  -0.82%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "54.36"
$ws.Range("E47").Value = "  +1.80%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0507"
$ws.Range("E48").Value = "  -0.11%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.406"
$ws.Range("E49").Value = "  -0.79%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₇0957"
$ws.Range("E50").Value = "  -8.53%  "
$ws.Range("E51").Value = "  -0.31%  "
